# Normalize the ordering of names/emails in the "Recorded By" column (G).
# Some cells list the recorder(s) in a different order than the canonical
# one used elsewhere in the report; this brings them in line by swapping
# specific known orderings to their corrected equivalent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of old (unordered) "Recorded By" text -> corrected ordering.
$map = @{
    "System, admin@admin.com" = "admin@admin.com, System";
    "dnasr281@gmail.com, System" = "System, dnasr281@gmail.com";
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com";
    "system, System, backup@backdoor.com" = "System, system, backup@backdoor.com";
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    if ($null -ne $current -and $map.ContainsKey($current)) {
        $cell.Value2 = $map[$current]
    }
}
